$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.183.30"
$ws.Range("E2").Value = "  +0.79%  "

$ws.Range("D3").Value = "1.843.54"
$ws.Range("E3").Value = "  +1.13%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9990"
$ws.Range("E4").Value = "  -0.55%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "279.73"
$ws.Range("E5").Value = "  -0.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9986"
$ws.Range("E6").Value = "  -0.52%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5096"
$ws.Range("E7").Value = "  +0.16%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3511"
$ws.Range("E8").Value = "  -1.58%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.93"
$ws.Range("E9").Value = "  -0.43%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06837"
$ws.Range("E10").Value = "  +1.85%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.08"
$ws.Range("E11").Value = "  -0.99%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.8072"
$ws.Range("E12").Value = "  -5.49%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07799"
$ws.Range("E13").Value = "  -1.23%  "

$ws.Range("D14").Value = "1.858.90"
$ws.Range("E14").Value = "  +1.21%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.094"
$ws.Range("E15").Value = "  +0.69%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.66"
$ws.Range("E16").Value = "  +0.98%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9988"
$ws.Range("E17").Value = "  -0.58%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.20"
$ws.Range("E18").Value = "  +0.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008070"
$ws.Range("E19").Value = "  -1.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9989"
$ws.Range("E20").Value = "  -0.60%  "

$ws.Range("D21").Value = "26.208.62"
$ws.Range("E21").Value = "  +0.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.775"
$ws.Range("E22").Value = "  +0.55%  "

$ws.Range("E23").Value = "  -0.76%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.219"
$ws.Range("E24").Value = "  +1.27%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.387"
$ws.Range("E25").Value = "  +10.92%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "144.29"
$ws.Range("E26").Value = "  +1.57%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.662"
$ws.Range("E27").Value = "  -1.48%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.23"
$ws.Range("E28").Value = "  +1.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "110.10"
$ws.Range("E29").Value = "  +0.82%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.372"
$ws.Range("E30").Value = "  +0.89%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.306"
$ws.Range("E31").Value = "  +1.38%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08748"
$ws.Range("E32").Value = "  -1.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04910"
$ws.Range("E33").Value = "  +2.33%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.172"
$ws.Range("E34").Value = "  +3.96%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7348"
$ws.Range("E35").Value = "  -1.38%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.851"
$ws.Range("E36").Value = "  +0.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.242"
$ws.Range("E37").Value = "  +4.77%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.399"
$ws.Range("E38").Value = "  -1.77%  "

$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5196"
$ws.Range("E39").Value = "  -4.83%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01851"
$ws.Range("E40").Value = "  -0.65%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9659"
$ws.Range("E41").Value = "  -1.86%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "116.17"
$ws.Range("E42").Value = "  +2.90%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.283"
$ws.Range("E43").Value = "  +1.09%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.024"
$ws.Range("E44").Value = "  -2.53%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9980"
$ws.Range("E45").Value = "  -0.64%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4525"
$ws.Range("E46").Value = "  -4.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1358"
$ws.Range("E47").Value = "  -1.82%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.327"
$ws.Range("E48").Value = "  +0.86%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.30"
$ws.Range("E49").Value = "  +1.40%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.503"
$ws.Range("E50").Value = "  +0.01%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05928"
$ws.Range("E51").Value = "  +0.16%  "
